$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents (ID, Classification, CMS?, URI, API Mapping (CL7))
$data = @(
    @("ID", "Classification", "CMS?", "URI", "API Mapping (CL7)"),
    @("SMKCO7001", "Decorative Art", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#DecorativeArt", "Kunsthåndværk"),
    @("SMKCO7002", "Drawing", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Drawing", "Tegning"),
    @("SMKCO7003", "Cast", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Cast", "Afstøbning"),
    @("SMKCO7004", "Media Art", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#MediaArt", "Audiovisuel kunst"),
    @("SMKCO7005", "Painting", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Painting", "Maleri"),
    @("SMKCO7006", "Photograph", "Y", " http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Photograph", "Fotografi"),
    @("SMKCO7007", "Installation", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Installation", "Installation"),
    @("SMKCO7008", "Print", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Print", "Grafik"),
    @("SMKCO7009", "Sculpture", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Sculpture", "Skulptur"),
    @("SMKCO7010", "Frame", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Frame", "Ramme"),
    @("SMKCO7011", "Works of Art on Paper in Multiple Parts", "Y", "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#WorksOfArtOnPaperInMultipleParts", "Works of Art on Paper in Multiple Parts"),
    @("SMKCO7012", "Collage", $null, "http://www.sirma.com/ontologies/2016/02/culturalHeritageConservation/SMK#Collage", "Collage")
)

# Remove the old row 14 entirely (table shrank from 14 to 13 rows)
$null = $ws.Rows.Item(14).Delete()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    for ($c = 0; $c -lt 5; $c++) {
        $val = $row[$c]
        if ($val -eq $null) {
            $ws.Cells.Item($r, $c + 1).Value = ""
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

$null = $ws.Range("D18").Select()
